# Add 14 new rows (597-610) of landscaping data to Sheet1, continuing the
# existing daily readings table for 2025-08-03 and 2025-08-04.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Row, A Date, B Plant_Type, C Plant_Size, D Low, E High, G Rain,
#          H Growth, I Pruned, J Quadrant, K Shade, L UV, M Humidity,
#          N Dew_Point, O Pressure, P Wind_Gust, Q Cloud_Cover,
#          R Visibility, S AQI, T Pollen
$data = @(
    @(597, 45872, "Flowering", "Large", 78, 91, 0, 0, "No", 2, "Bright", 9, 0.56999999999999995, 79, 30.05, 12, 0.47, 8.1, 45, 4),
    @(598, 45872, "Nonflowering", "Medium", 78, 91, 0, 0, "No", 3, "Bright", 9, 0.56999999999999995, 79, 30.05, 12, 0.47, 8.1, 45, 4),
    @(599, 45872, "Nonflowering", "Small", 78, 91, 0, 0, "No", 3, "Neutral", 9, 0.56999999999999995, 79, 30.05, 12, 0.47, 8.1, 45, 4),
    @(600, 45872, "Nonflowering", "Medium", 78, 91, 0, 0, "No", 3, "Neutral", 9, 0.56999999999999995, 79, 30.05, 12, 0.47, 8.1, 45, 4),
    @(601, 45872, "Nonflowering", "Medium", 78, 91, 0, 0, "No", 3, "Dark", 9, 0.56999999999999995, 79, 30.05, 12, 0.47, 8.1, 45, 4),
    @(602, 45872, "Nonflowering", "Large", 78, 91, 0, 0, "No", 4, "Bright", 9, 0.56999999999999995, 79, 30.05, 12, 0.47, 8.1, 45, 4),
    @(603, 45872, "Tree", "Medium", 78, 91, 0, 0, "No", 1, "Bright", 9, 0.56999999999999995, 79, 30.05, 12, 0.47, 8.1, 45, 4),
    @(604, 45873, "Flowering", "Large", 79, 90, 0, 0, "No", 2, "Neutral", 10, 0.72, 79, 30, 15, 0.17, 8.1, 42, 1),
    @(605, 45873, "Nonflowering", "Medium", 79, 90, 0, 0, "No", 3, "Neutral", 10, 0.72, 79, 30, 15, 0.17, 8.1, 42, 1),
    @(606, 45873, "Nonflowering", "Small", 79, 90, 0, 0, "No", 3, "Bright", 10, 0.72, 79, 30, 15, 0.17, 8.1, 42, 1),
    @(607, 45873, "Nonflowering", "Medium", 79, 90, 0, 0, "No", 3, "Neutral", 10, 0.72, 79, 30, 15, 0.17, 8.1, 42, 1),
    @(608, 45873, "Nonflowering", "Medium", 79, 90, 0, 0, "No", 3, "Bright", 10, 0.72, 79, 30, 15, 0.17, 8.1, 42, 1),
    @(609, 45873, "Nonflowering", "Large", 79, 90, 0, 0, "No", 4, "Dark", 10, 0.72, 79, 30, 15, 0.17, 8.1, 42, 1),
    @(610, 45873, "Tree", "Medium", 79, 90, 0, 0, "No", 1, "Bright", 10, 0.72, 79, 30, 15, 0.17, 8.1, 42, 1)
)

foreach ($row in $data) {
    $r = $row[0]

    # Date column keeps the same date-number-format style as the rows above it.
    $ws.Cells.Item(596, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $row[1]

    $ws.Cells.Item($r, 2).Value  = $row[2]
    $ws.Cells.Item($r, 3).Value  = $row[3]
    $ws.Cells.Item($r, 4).Value  = $row[4]
    $ws.Cells.Item($r, 5).Value  = $row[5]
    $ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
    $ws.Cells.Item($r, 7).Value  = $row[6]
    $ws.Cells.Item($r, 8).Value  = $row[7]
    $ws.Cells.Item($r, 9).Value  = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
    $ws.Cells.Item($r, 14).Value = $row[13]
    $ws.Cells.Item($r, 15).Value = $row[14]
    $ws.Cells.Item($r, 16).Value = $row[15]
    $ws.Cells.Item($r, 17).Value = $row[16]
    $ws.Cells.Item($r, 18).Value = $row[17]
    $ws.Cells.Item($r, 19).Value = $row[18]
    $ws.Cells.Item($r, 20).Value = $row[19]
}

# Scroll / selection state, matching what Excel shows after adding the rows.
$ws.Application.ActiveWindow.ScrollRow = 584
$ws.Range("P604:P610").Select()

Write-Host "Added rows 597-610"
